# Update the crypto price/volume table (columns D = Price, E = Volume(1h))
# with the latest scraped figures, per the automated GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'28.252.44"
$ws.Cells.Item(2, 5).Value = "  +1.35%  "
$ws.Cells.Item(3, 4).Value = "'1.806.20"
$ws.Cells.Item(3, 5).Value = "  +3.14%  "
$ws.Cells.Item(4, 4).Value = "'1.003"
$ws.Cells.Item(4, 5).Value = "  -0.21%  "
$ws.Cells.Item(5, 4).Value = "'336.58"
$ws.Cells.Item(5, 5).Value = "  +0.42%  "
$ws.Cells.Item(6, 4).Value = "'0.9994"
$ws.Cells.Item(6, 5).Value = "  -0.10%  "
$ws.Cells.Item(7, 4).Value = "'0.4618"
$ws.Cells.Item(7, 5).Value = "  +20.69%  "
$ws.Cells.Item(8, 4).Value = "'0.3726"
$ws.Cells.Item(8, 5).Value = "  +9.52%  "
$ws.Cells.Item(9, 4).Value = "'45.15"
$ws.Cells.Item(9, 5).Value = "  -2.57%  "
$ws.Cells.Item(10, 4).Value = "'0.07669"
$ws.Cells.Item(10, 5).Value = "  +6.24%  "
$ws.Cells.Item(11, 4).Value = "'1.152"
$ws.Cells.Item(11, 5).Value = "  +3.47%  "
$ws.Cells.Item(12, 4).Value = "'22.39"
$ws.Cells.Item(12, 5).Value = "  -0.28%  "
$ws.Cells.Item(13, 5).Value = "  -0.20%  "
$ws.Cells.Item(14, 4).Value = "'6.348"
$ws.Cells.Item(14, 5).Value = "  +3.07%  "
$ws.Cells.Item(15, 4).Value = "'7.487"
$ws.Cells.Item(15, 5).Value = "  +4.91%  "
$ws.Cells.Item(16, 4).Value = "'1.806.14"
$ws.Cells.Item(16, 5).Value = "  +2.93%  "
$ws.Cells.Item(17, 4).Value = "'0.00001099"
$ws.Cells.Item(17, 5).Value = "  +3.72%  "
$ws.Cells.Item(18, 4).Value = "'0.06724"
$ws.Cells.Item(18, 5).Value = "  +1.77%  "
$ws.Cells.Item(19, 4).Value = "'82.00"
$ws.Cells.Item(19, 5).Value = "  +4.04%  "
$ws.Cells.Item(20, 4).Value = "'0.9992"
$ws.Cells.Item(20, 5).Value = "  -0.16%  "
$ws.Cells.Item(21, 4).Value = "'17.47"
$ws.Cells.Item(21, 5).Value = "  +4.69%  "
$ws.Cells.Item(22, 4).Value = "'6.423"
$ws.Cells.Item(22, 5).Value = "  +3.23%  "
$ws.Cells.Item(23, 4).Value = "'28.244.58"
$ws.Cells.Item(23, 5).Value = "  +1.24%  "
$ws.Cells.Item(24, 5).Value = "  +1.93%  "
$ws.Cells.Item(25, 4).Value = "'2.411"
$ws.Cells.Item(25, 5).Value = "  +1.23%  "
$ws.Cells.Item(26, 4).Value = "'20.91"
$ws.Cells.Item(26, 5).Value = "  +5.43%  "
$ws.Cells.Item(27, 4).Value = "'153.93"
$ws.Cells.Item(27, 5).Value = "  +0.51%  "
$ws.Cells.Item(28, 5).Value = "  +3.03%  "
$ws.Cells.Item(29, 4).Value = "'2.012.04"
$ws.Cells.Item(29, 5).Value = "  +2.93%  "
$ws.Cells.Item(30, 4).Value = "'133.51"
$ws.Cells.Item(30, 5).Value = "  +1.05%  "
$ws.Cells.Item(31, 4).Value = "'1.259"
$ws.Cells.Item(31, 5).Value = "  -0.62%  "
$ws.Cells.Item(32, 4).Value = "'4.033"
$ws.Cells.Item(32, 5).Value = "  +0.31%  "
$ws.Cells.Item(33, 4).Value = "'0.09567"
$ws.Cells.Item(33, 5).Value = "  +8.64%  "
$ws.Cells.Item(34, 4).Value = "'5.872"
$ws.Cells.Item(34, 5).Value = "  +0.71%  "
$ws.Cells.Item(35, 4).Value = "'0.2223"
$ws.Cells.Item(35, 5).Value = "  +5.85%  "
$ws.Cells.Item(36, 4).Value = "'12.14"
$ws.Cells.Item(36, 5).Value = "  -0.46%  "
$ws.Cells.Item(37, 4).Value = "'0.06369"
$ws.Cells.Item(37, 5).Value = "  +3.57%  "
$ws.Cells.Item(38, 4).Value = "'0.02356"
$ws.Cells.Item(38, 5).Value = "  +3.12%  "
$ws.Cells.Item(39, 4).Value = "'5.265"
$ws.Cells.Item(39, 5).Value = "  +2.44%  "
$ws.Cells.Item(40, 4).Value = "'0.6655"
$ws.Cells.Item(40, 5).Value = "  +1.33%  "
$ws.Cells.Item(41, 5).Value = "  +0.52%  "
$ws.Cells.Item(42, 4).Value = "'1.239"
$ws.Cells.Item(42, 5).Value = "  +2.50%  "
$ws.Cells.Item(43, 4).Value = "'8.271"
$ws.Cells.Item(43, 5).Value = "  +3.50%  "
$ws.Cells.Item(44, 4).Value = "'14.39"
$ws.Cells.Item(44, 5).Value = "  +4.99%  "
$ws.Cells.Item(45, 4).Value = "'0.9989"
$ws.Cells.Item(45, 5).Value = "  -0.17%  "
$ws.Cells.Item(46, 4).Value = "'0.6129"
$ws.Cells.Item(46, 5).Value = "  +1.03%  "
$ws.Cells.Item(47, 4).Value = "'3.826"
$ws.Cells.Item(47, 5).Value = "  +0.13%  "
$ws.Cells.Item(48, 4).Value = "'129.86"
$ws.Cells.Item(48, 5).Value = "  +2.69%  "
$ws.Cells.Item(49, 4).Value = "'2.053"
$ws.Cells.Item(49, 5).Value = "  +2.60%  "
$ws.Cells.Item(50, 4).Value = "'0.07163"
$ws.Cells.Item(50, 5).Value = "  +2.69%  "
$ws.Cells.Item(51, 4).Value = "'1.180"
$ws.Cells.Item(51, 5).Value = "  +0.75%  "
